# Auto-generated edit script: apply numeric updates to Shinryu_Profits workbook
# Source: diff of Sheets/Shinryu_Profits.xlsx (per-sheet leve profit recalculation)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 66666828
$ws.Range("I41").Value = 111111190
$ws.Range("J41").Value = 276
$ws.Range("K41").Value = 111111190
$ws.Range("L41").Value = 276
$ws.Range("M41").Value = -111110750
$ws.Range("N41").Value = -1156

$ws.Range("H51").Value = 1471.2858
$ws.Range("I51").Value = 1666.3334
$ws.Range("J51").Value = 1325
$ws.Range("K51").Value = 1666.3334
$ws.Range("L51").Value = 1325
$ws.Range("M51").Value = -1182.3334
$ws.Range("N51").Value = -2293

$ws.Range("H63").Value = 35000
$ws.Range("J63").Value = 35000
$ws.Range("L63").Value = 35000
$ws.Range("N63").Value = -36248

$ws.Range("H66").Value = 35000
$ws.Range("J66").Value = 35000
$ws.Range("L66").Value = 105000
$ws.Range("N66").Value = -111240

$ws.Range("H75").Value = 19618.2
$ws.Range("J75").Value = 19618.2
$ws.Range("L75").Value = 19618.2
$ws.Range("N75").Value = -21490.2

$ws.Range("H76").Value = 3499.2104
$ws.Range("I76").Value = 2998.5
$ws.Range("J76").Value = 3949.85
$ws.Range("K76").Value = 2998.5
$ws.Range("L76").Value = 3949.85
$ws.Range("M76").Value = -2683.5
$ws.Range("N76").Value = -4579.85

$ws.Range("H78").Value = 19618.2
$ws.Range("J78").Value = 19618.2
$ws.Range("L78").Value = 58854.60000000001
$ws.Range("N78").Value = -68214.60000000001

$ws.Range("H79").Value = 3499.2104
$ws.Range("I79").Value = 2998.5
$ws.Range("J79").Value = 3949.85
$ws.Range("K79").Value = 2998.5
$ws.Range("L79").Value = 3949.85
$ws.Range("M79").Value = -1906.5
$ws.Range("N79").Value = -6133.85

$ws.Range("H86").Value = 1788.1111
$ws.Range("I86").Value = 1786.4
$ws.Range("J86").Value = 1796.6666
$ws.Range("K86").Value = 1786.4
$ws.Range("L86").Value = 1796.6666
$ws.Range("M86").Value = -663.4000000000001
$ws.Range("N86").Value = -4042.6666

$ws.Range("H89").Value = 1788.1111
$ws.Range("I89").Value = 1786.4
$ws.Range("J89").Value = 1796.6666
$ws.Range("K89").Value = 8932
$ws.Range("L89").Value = 8983.333000000001
$ws.Range("M89").Value = -3316
$ws.Range("N89").Value = -20215.333

$ws.Range("H92").Value = 1282.2858
$ws.Range("I92").Value = 412.66666
$ws.Range("J92").Value = 6500
$ws.Range("K92").Value = 412.66666
$ws.Range("L92").Value = 6500
$ws.Range("M92").Value = 835.33334
$ws.Range("N92").Value = -8996

$ws.Range("H125").Value = 1636
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H127").Value = 869.4666999999999
$ws.Range("I127").Value = 389.5
$ws.Range("J127").Value = 1189.4445
$ws.Range("K127").Value = 1168.5
$ws.Range("L127").Value = 3568.3335
$ws.Range("M127").Value = 3791.5
$ws.Range("N127").Value = -13488.3335

$ws.Range("H131").Value = 26043.414
$ws.Range("I131").Value = 32125.969
$ws.Range("J131").Value = 4416.5557
$ws.Range("K131").Value = 96377.90700000001
$ws.Range("L131").Value = 13249.6671
$ws.Range("M131").Value = -91337.90700000001
$ws.Range("N131").Value = -23329.6671

$ws.Range("H132").Value = 2949.7083
$ws.Range("I132").Value = 2989.35
$ws.Range("J132").Value = 2751.5
$ws.Range("K132").Value = 8968.049999999999
$ws.Range("L132").Value = 8254.5
$ws.Range("M132").Value = -6438.049999999999
$ws.Range("N132").Value = -13314.5

$ws.Range("H138").Value = 3777.7834
$ws.Range("I138").Value = 808.7
$ws.Range("J138").Value = 5262.325
$ws.Range("K138").Value = 2426.1
$ws.Range("L138").Value = 15786.975
$ws.Range("M138").Value = 2713.9
$ws.Range("N138").Value = -26066.975


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18522792
$ws.Range("I32").Value = 20410236
$ws.Range("J32").Value = 25840
$ws.Range("K32").Value = 20410236
$ws.Range("L32").Value = 25840
$ws.Range("M32").Value = -20409949
$ws.Range("N32").Value = -26414

$ws.Range("H88").Value = 2112.889
$ws.Range("I88").Value = 2219.3333
$ws.Range("J88").Value = 1900
$ws.Range("K88").Value = 2219.3333
$ws.Range("L88").Value = 1900
$ws.Range("M88").Value = -1813.3333
$ws.Range("N88").Value = -2712

$ws.Range("H91").Value = 2112.889
$ws.Range("I91").Value = 2219.3333
$ws.Range("J91").Value = 1900
$ws.Range("K91").Value = 2219.3333
$ws.Range("L91").Value = 1900
$ws.Range("M91").Value = -815.3332999999998
$ws.Range("N91").Value = -4708

$ws.Range("H97").Value = 627.36365
$ws.Range("I97").Value = 640.1
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 640.1
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = -144.1
$ws.Range("N97").Value = -1492


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2815
$ws.Range("I86").Value = 2412.5
$ws.Range("J86").Value = 3083.3333
$ws.Range("K86").Value = 2412.5
$ws.Range("L86").Value = 3083.3333
$ws.Range("M86").Value = -1289.5
$ws.Range("N86").Value = -5329.3333

$ws.Range("H89").Value = 2815
$ws.Range("I89").Value = 2412.5
$ws.Range("J89").Value = 3083.3333
$ws.Range("K89").Value = 12062.5
$ws.Range("L89").Value = 15416.6665
$ws.Range("M89").Value = -6446.5
$ws.Range("N89").Value = -26648.6665


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2214.9365
$ws.Range("I31").Value = 1643.4906
$ws.Range("K31").Value = 1643.4906
$ws.Range("M31").Value = -1348.4906

$ws.Range("H34").Value = 2214.9365
$ws.Range("I34").Value = 1643.4906
$ws.Range("K34").Value = 1643.4906
$ws.Range("M34").Value = -1441.4906


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1101.4
$ws.Range("J86").Value = 1176.25
$ws.Range("L86").Value = 3528.75
$ws.Range("N86").Value = -5900.75

$ws.Range("H89").Value = 1101.4
$ws.Range("J89").Value = 1176.25
$ws.Range("L89").Value = 10586.25
$ws.Range("N89").Value = -22442.25

$ws.Range("H98").Value = 206.58333
$ws.Range("I98").Value = 237.8
$ws.Range("J98").Value = 184.28572
$ws.Range("K98").Value = 713.4000000000001
$ws.Range("L98").Value = 552.85716
$ws.Range("M98").Value = 784.5999999999999
$ws.Range("N98").Value = -3548.85716

$ws.Range("H131").Value = 2413.492
$ws.Range("I131").Value = 517.6923
$ws.Range("J131").Value = 2906.4
$ws.Range("K131").Value = 1553.0769
$ws.Range("L131").Value = 8719.200000000001
$ws.Range("M131").Value = 3486.9231
$ws.Range("N131").Value = -18799.2


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3245.5
$ws.Range("I68").Value = 2501
$ws.Range("J68").Value = 3990
$ws.Range("K68").Value = 2501
$ws.Range("L68").Value = 3990
$ws.Range("M68").Value = -1752
$ws.Range("N68").Value = -5488

$ws.Range("H71").Value = 3245.5
$ws.Range("I71").Value = 2501
$ws.Range("J71").Value = 3990
$ws.Range("K71").Value = 12505
$ws.Range("L71").Value = 19950
$ws.Range("M71").Value = -8761
$ws.Range("N71").Value = -27438

$ws.Range("H93").Value = 7514.875
$ws.Range("I93").Value = 7949.2
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 7949.2
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -6701.2
$ws.Range("N93").Value = -3496


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11891.454
$ws.Range("I62").Value = 21920.6
$ws.Range("J62").Value = 3533.8333
$ws.Range("K62").Value = 21920.6
$ws.Range("L62").Value = 3533.8333
$ws.Range("M62").Value = -21296.6
$ws.Range("N62").Value = -4781.8333

$ws.Range("H65").Value = 11891.454
$ws.Range("I65").Value = 21920.6
$ws.Range("J65").Value = 3533.8333
$ws.Range("K65").Value = 109603
$ws.Range("L65").Value = 17669.1665
$ws.Range("M65").Value = -106483
$ws.Range("N65").Value = -23909.1665

$ws.Range("H96").Value = 2336.2727
$ws.Range("I96").Value = 1559.8
$ws.Range("J96").Value = 2983.3333
$ws.Range("K96").Value = 1559.8
$ws.Range("L96").Value = 2983.3333
$ws.Range("M96").Value = -186.8
$ws.Range("N96").Value = -5729.3333

